$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edit: "Department" column header becomes "Department Code" ---
$ws.Range("F1").Value = "Department Code"

# --- Content edit: new entry "CRR" added as a second data row ---
$ws.Range("A2").Value = "CRR"

# --- Header row formatting: no longer bold, bumped up to 12pt ---
$hdr = $ws.Range("A1:I1")
$hdr.Font.Bold = $false
$hdr.Font.Size = 12
$ws.Range("A2").Font.Bold = $false
$ws.Range("A2").Font.Size = 12
